$d = $word.ActiveDocument

# Locate the paragraph that contains the "LOB1012" requirement text, then
# remove the following three paragraphs:
#   1) the blank paragraph right after it
#   2) the blank paragraph that carries the page-break-before
#   3) the "(c) 2020 ... Creative Commons Attribution" paragraph
$count = $d.Paragraphs.Count
for ($i = $count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "*LOB1012*") {
        $target = $d.Range($p.Range.End, $d.Paragraphs.Item($i + 3).Range.End)
        $target.Delete()
        break
    }
}
